$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = -2.0186232039417233
$ws.Range("C3").Value = -1.9224434696129178
$ws.Range("B4").Value = -2.1005791643768208
$ws.Range("C4").Value = -1.8775777800155538
$ws.Range("B5").Value = -2.2303034310508951
$ws.Range("C5").Value = -1.9121306647320995
$ws.Range("B6").Value = -2.3114680809384596
$ws.Range("C6").Value = -1.916388542554772
$ws.Range("B7").Value = -2.1556543707971789
$ws.Range("C7").Value = -1.6877873227614766
$ws.Range("B9").Value = 0.20147957314176371
$ws.Range("C9").Value = 0.16404870315790221
$ws.Range("B13").Value = 0.046394883432239062
$ws.Range("C13").Value = 0.093926923652357941
$ws.Range("B14").Value = 0.85475969645446448
$ws.Range("C14").Value = 0.98521409797073434
$ws.Range("C16").Value = -0.095021447472647583
$ws.Range("C17").Value = 1.1464536799357905
$ws.Range("C18").Value = 0.46918149865077119
$ws.Range("C19").Value = -1.5150641726209138
$ws.Range("C20").Value = 0.0099225211219549636
$ws.Range("B21").Value = 8043
$ws.Range("C21").Value = 9990
